$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.362.95"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.842.08"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'259.16"
$ws.Range("E5").Value = "  -7.62%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("D7").Value = "'0.5097"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").Value = "'0.3199"
$ws.Range("E8").Value = "  -8.86%  "
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "'18.94"
$ws.Range("E10").Value = "  -5.31%  "
$ws.Range("D11").Value = "'0.7700"
$ws.Range("E11").Value = "  -5.19%  "
$ws.Range("D12").Value = "'0.07686"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "1.881.83"
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "'88.27"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "'5.017"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").Value = "'14.04"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "26.392.24"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "2.076.60"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'4.580"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("D23").Value = "'9.533"
$ws.Range("E23").Value = "  -5.69%  "
$ws.Range("D24").Value = "'5.979"
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("D25").Value = "'2.331"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("D26").Value = "'145.18"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'1.651"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'16.91"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "'110.83"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("D30").Value = "'4.167"
$ws.Range("E30").Value = "  -4.74%  "
$ws.Range("D31").Value = "'4.162"
$ws.Range("E31").Value = "  -3.57%  "
$ws.Range("D32").Value = "'0.08702"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "'0.04814"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").Value = "'1.130"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("D35").Value = "'2.840"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "'0.6859"
$ws.Range("E36").Value = "  -7.69%  "
$ws.Range("D37").Value = "'3.084"
$ws.Range("E37").Value = "  -5.09%  "
$ws.Range("D38").Value = "'0.01806"
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").Value = "'2.207"
$ws.Range("E39").Value = "  -7.70%  "
$ws.Range("D40").Value = "'0.4892"
$ws.Range("E40").Value = "  -5.57%  "
$ws.Range("D41").Value = "'113.16"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").Value = "'0.9043"
$ws.Range("E42").Value = "  -7.01%  "
$ws.Range("D43").Value = "'6.121"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'7.753"
$ws.Range("E45").Value = "  -3.59%  "
$ws.Range("D46").Value = "'0.4240"
$ws.Range("E46").Value = "  -6.52%  "
$ws.Range("D47").Value = "'0.1268"
$ws.Range("E47").Value = "  -6.79%  "
$ws.Range("D48").Value = "'9.157"
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("D49").Value = "'0.05891"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "'35.08"
$ws.Range("E50").Value = "  -3.40%  "
$ws.Range("D51").Value = "'1.425"
$ws.Range("E51").Value = "  -5.27%  "
